# "cambios de agosto, puntos fe de ratas e historico"
# Update the "Fecha de inicio/termino del periodo", "Fecha de elaboracion",
# "Fecha de validacion" and "Fecha de actualizacion" columns on the single
# data row (row 8) of the report, and move the viewport/selection from
# D12 over to G13 (columns G:I, "Hipervinculo/Area responsable" section).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fecha de inicio del periodo que se informa
$ws.Range("B8").Value = 44652
# Fecha de termino del periodo que se informa
$ws.Range("C8").Value = 44742
# Fecha de elaboracion
$ws.Range("E8").Value = 44742
# Fecha de validacion
$ws.Range("H8").Value = 44753
# Fecha de actualizacion
$ws.Range("I8").Value = 44753

# Scroll the window so column G is left-most (was column D) ...
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 2
# ... and move the active selection to G13 (was D12).
$ws.Range("G13").Select()
